# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets,
# reflecting refreshed figures from the upstream data source.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1037
$ws1.Range("F4").Value = 174
$ws1.Range("F5").Value = 2821
$ws1.Range("F10").Value = 80
$ws1.Range("F12").Value = 2653
$ws1.Range("F13").Value = 847

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1037
$ws4.Range("F5").Value = 174
$ws4.Range("F6").Value = 2821
$ws4.Range("F12").Value = 80
$ws4.Range("F14").Value = 2653
$ws4.Range("F15").Value = 847
